# Insert two new data rows right before the current row 249 ("Femacal de La
# Calera" / Zapallo sheet), pushing the existing rows 249-356 down to 251-358.
# This mirrors a weekly data refresh: two fresh observations (Camote / Paine,
# dated 44510) are prepended to the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 249 (each Insert() pushes the row that is
# currently at 249, and everything below it, down by one).
$ws.Rows.Item(249).Insert()
$ws.Rows.Item(250).Insert()

# Copy the date-number formatting from the row that now sits just below the
# freshly inserted rows (the old row 249, now at 251) onto the new D cells so
# they keep the same number format (style index 2 in the original file).
$ws.Range("D249").NumberFormat = $ws.Range("D251").NumberFormat
$ws.Range("D250").NumberFormat = $ws.Range("D251").NumberFormat

# --- New row 249: Camote, 1a nueva(o) ---
$ws.Range("A249").Value = 3
$ws.Range("B249").Value = "Femacal de La Calera"
$ws.Range("C249").Value = "Coquimbo"
$ws.Range("D249").Value = 44510
$ws.Range("E249").Value = 5
$ws.Range("F249").Value = 100112045
$ws.Range("G249").Value = "Zapallo"
$ws.Range("H249").Value = "Camote"
$ws.Range("I249").Value = "1a nueva(o)"
$ws.Range("J249").Value = 165
$ws.Range("K249").Value = 500
$ws.Range("L249").Value = 550
$ws.Range("M249").Value = 526
$ws.Range("N249").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O249").Value = "Perú"
$ws.Range("P249").Value = 526
$ws.Range("Q249").Value = 1
$ws.Range("R249").Value = "Hortaliza"

# --- New row 250: Paine, 1a nueva(o) ---
$ws.Range("A250").Value = 3
$ws.Range("B250").Value = "Femacal de La Calera"
$ws.Range("C250").Value = "Coquimbo"
$ws.Range("D250").Value = 44510
$ws.Range("E250").Value = 5
$ws.Range("F250").Value = 100112045
$ws.Range("G250").Value = "Zapallo"
$ws.Range("H250").Value = "Paine"
$ws.Range("I250").Value = "1a nueva(o)"
$ws.Range("J250").Value = 110
$ws.Range("K250").Value = 150
$ws.Range("L250").Value = 150
$ws.Range("M250").Value = 150
$ws.Range("N250").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O250").Value = "Provincia de Talca"
$ws.Range("P250").Value = 150
$ws.Range("Q250").Value = 1
$ws.Range("R250").Value = "Hortaliza"
